$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.164.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.69%  '
$ws.Range("D3").Value = '''2.253.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.78%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''244.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("D6").Value = '''0.617'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("D7").Value = '''76.22'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +9.86%  '
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '''0.612'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.75%  '
$ws.Range("D10").Value = '''41.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.92%  '
$ws.Range("D11").Value = '''0.0938'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.44%  '
$ws.Range("D12").Value = '''7.02'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.69%  '
$ws.Range("D14").Value = '''2.591.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.71%  '
$ws.Range("D15").Value = '''14.65'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.01%  '
$ws.Range("D16").Value = '''2.245.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.83%  '
$ws.Range("D17").Value = '''0.806'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.69%  '
$ws.Range("D18").Value = '''43.065.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.02%  '
$ws.Range("D19").Value = '''0.0000105'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.52%  '
$ws.Range("D20").Value = '''71.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("D21").Value = '''6.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.25%  '
$ws.Range("D22").Value = '''10.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.16%  '
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").Value = '''2.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +16.55%  '
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '''230.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").Value = '''10.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.37%  '
$ws.Range("D27").Value = '''3.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.18%  '
$ws.Range("D28").Value = '''39.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +32.44%  '
$ws.Range("E29").Value = '  +2.62%  '
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("D31").Value = '''174.05'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.90%  '
$ws.Range("D32").Value = '''20.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.82%  '
$ws.Range("E33").Value = '  +5.06%  '
$ws.Range("D34").Value = '''5.38'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.74%  '
$ws.Range("D35").Value = '''0.123'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.86%  '
$ws.Range("D36").Value = '''0.110'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.67%  '
$ws.Range("D38").Value = '''0.0333'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +18.73%  '
$ws.Range("D39").Value = '''13.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +13.51%  '
$ws.Range("E40").Value = '  +4.38%  '
$ws.Range("E41").Value = '  +3.16%  '
$ws.Range("D42").Value = '''0.204'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.09%  '
$ws.Range("D43").Value = '''60.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.36%  '
$ws.Range("D44").Value = '''106.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.51%  '
$ws.Range("D45").Value = '''8.72'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.76%  '
$ws.Range("D46").Value = '''0.100'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.46%  '
$ws.Range("D47").Value = '''0.476'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +29.78%  '
$ws.Range("D48").Value = '''2.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.98%  '
$ws.Range("E49").Value = '  +3.61%  '
$ws.Range("E50").Value = '  +2.65%  '
$ws.Range("D51").Value = '''2.463.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.76%  '
